# Update countries & provincias Spain
#
# Refreshes the COVID-19 "Pais" table with a newer data pull:
#   - bumps the "last updated" timestamp in A1
#   - updates totals/new-cases/active/recovered/critical/deaths-today/deaths
#     for the countries whose figures moved
#   - three countries (Kuwait/Belgica, Bosnia y Herzegovina/Corea del Sur,
#     Libia/Republica de Macedonia) leapfrog their neighbour in the ranking,
#     so the two rows' country names trade places along with their numbers

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (A1) ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 15:15"

# --- Row 4: Estados Unidos ---------------------------------------------
$ws.Range("B4").Value = 6212708
$ws.Range("C4").Value = 912
$ws.Range("D4").Value = 3457318
$ws.Range("E4").Value = 2567597
$ws.Range("G4").Value = 57
$ws.Range("H4").Value = 187793

# --- Row 13: Argentina --------------------------------------------------
$ws.Range("D13").Value = 308376
$ws.Range("E13").Value = 100629
$ws.Range("G13").Value = 70
$ws.Range("H13").Value = 8730

# --- Row 17: Arabia Saudita ---------------------------------------------
$ws.Range("B17").Value = 316670
$ws.Range("C17").Value = 898
$ws.Range("D17").Value = 291514
$ws.Range("E17").Value = 21227
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 3929

# --- Row 24: Irak --------------------------------------------------------
$ws.Range("B24").Value = 238338
$ws.Range("C24").Value = 3404
$ws.Range("D24").Value = 180473
$ws.Range("E24").Value = 50742
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = 7123

# --- Rows 39/40: Kuwait overtakes Belgica --------------------------------
$ws.Range("A39").Value = "Kuwait"
$ws.Range("B39").Value = 85811
$ws.Range("C39").Value = 702
$ws.Range("D39").Value = 77657
$ws.Range("E39").Value = 7620
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 534

$ws.Range("A40").Value = "Belgica"
$ws.Range("B40").Value = 85236
$ws.Range("C40").Value = 194
$ws.Range("D40").Value = 18422
$ws.Range("E40").Value = 56919
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 9895

# --- Row 42: Suecia -------------------------------------------------------
$ws.Range("B42").Value = 84521
$ws.Range("G42").Value = 5
$ws.Range("H42").Value = 5813

# --- Row 45: Paises Bajos --------------------------------------------------
$ws.Range("B45").Value = 71129
$ws.Range("C45").Value = 462
$ws.Range("G45").Value = 6
$ws.Range("H45").Value = 6230

# --- Row 61: Suiza -----------------------------------------------------------
$ws.Range("E61").Value = 4282
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 2011

# --- Row 67: Azerbaiyan -------------------------------------------------------
$ws.Range("B67").Value = 36578
$ws.Range("C67").Value = 143
$ws.Range("D67").Value = 33977
$ws.Range("E67").Value = 2065
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 536

# --- Row 69: Serbia ------------------------------------------------------------
$ws.Range("B69").Value = 31482
$ws.Range("C69").Value = 76
$ws.Range("D69").Value = 30053
$ws.Range("E69").Value = 714
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 715

# --- Rows 76/77: Bosnia y Herzegovina overtakes Corea del Sur ------------------
$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("B76").Value = 20234
$ws.Range("C76").Value = 270
$ws.Range("D76").Value = 13435
$ws.Range("E76").Value = 6179
$ws.Range("G76").Value = 11
$ws.Range("H76").Value = 620

$ws.Range("A77").Value = "Corea del Sur"
$ws.Range("B77").Value = 20182
$ws.Range("C77").Value = 235
$ws.Range("D77").Value = 15198
$ws.Range("E77").Value = 4660
$ws.Range("H77").Value = 324

# --- Row 82: Dinamarca -----------------------------------------------------------
$ws.Range("B82").Value = 17084
$ws.Range("C82").Value = 99
$ws.Range("D82").Value = 15300
$ws.Range("E82").Value = 1159
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 625

# --- Rows 85/86: Libia overtakes Republica de Macedonia --------------------------
$ws.Range("A85").Value = "Libia"
$ws.Range("B85").Value = 14624
$ws.Range("C85").Value = 658
$ws.Range("D85").Value = 1676
$ws.Range("E85").Value = 12706
$ws.Range("G85").Value = 5
$ws.Range("H85").Value = 242

$ws.Range("A86").Value = "Republica de Macedonia"
$ws.Range("B86").Value = 14455
$ws.Range("C86").Value = 114
$ws.Range("D86").Value = 11525
$ws.Range("E86").Value = 2326
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 604

# --- Row 141: Islandia -------------------------------------------------------------
$ws.Range("B141").Value = 2116
$ws.Range("C141").Value = 9
$ws.Range("D141").Value = 2007
$ws.Range("E141").Value = 99

# --- Row 163: Vietnam ---------------------------------------------------------------
$ws.Range("D163").Value = 735
$ws.Range("E163").Value = 275
